$wb = $excel.ActiveWorkbook

# --- "summary" sheet (sheet1) numeric corrections ---
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("D2").Value = 1985
$wsSummary.Range("E2").Value = 81.21
$wsSummary.Range("F2").Value = -1.86
$wsSummary.Range("J2").Value = 0.24
$wsSummary.Range("L2").Value = 1.09
$wsSummary.Range("D3").Value = 1953
$wsSummary.Range("I3").Value = -1.74
$wsSummary.Range("J3").Value = 0.26
$wsSummary.Range("L3").Value = 1.05
$wsSummary.Range("D4").Value = 1862
$wsSummary.Range("E4").Value = 49.73
$wsSummary.Range("F4").Value = 0.03
$wsSummary.Range("I4").Value = 1.09
$wsSummary.Range("J4").Value = 0.26
$wsSummary.Range("L4").Value = 1.07
$wsSummary.Range("D5").Value = 1716
$wsSummary.Range("F5").Value = -0.01
$wsSummary.Range("H5").Value = 0.95
$wsSummary.Range("I5").Value = -2.34
$wsSummary.Range("J5").Value = 0.28
$wsSummary.Range("L5").Value = 0.86
$wsSummary.Range("D6").Value = 1379
$wsSummary.Range("E6").Value = 18.56
$wsSummary.Range("F6").Value = 1.89
$wsSummary.Range("H6").Value = 1.05
$wsSummary.Range("I6").Value = 1.11
$wsSummary.Range("L6").Value = 0.94
$wsSummary.Range("C7").Value = 680
$wsSummary.Range("D7").Value = 657
$wsSummary.Range("E7").Value = 75.95
$wsSummary.Range("F7").Value = -1.68
$wsSummary.Range("H7").Value = 1.04
$wsSummary.Range("I7").Value = 0.8
$wsSummary.Range("J7").Value = 0.26
$wsSummary.Range("K7").Value = 0.04
$wsSummary.Range("L7").Value = 0.95
$wsSummary.Range("C8").Value = 680
$wsSummary.Range("D8").Value = 655
$wsSummary.Range("E8").Value = 72.52
$wsSummary.Range("I8").Value = 1.7
$wsSummary.Range("J8").Value = 0.23
$wsSummary.Range("C9").Value = 680
$wsSummary.Range("D9").Value = 647
$wsSummary.Range("E9").Value = 64.91
$wsSummary.Range("F9").Value = -1.03
$wsSummary.Range("I9").Value = 0.42
$wsSummary.Range("J9").Value = 0.32
$wsSummary.Range("L9").Value = 1.12
$wsSummary.Range("C10").Value = 680
$wsSummary.Range("D10").Value = 648
$wsSummary.Range("E10").Value = 62.35
$wsSummary.Range("F10").Value = -0.88
$wsSummary.Range("I10").Value = 1.17
$wsSummary.Range("J10").Value = 0.28
$wsSummary.Range("L10").Value = 0.97
$wsSummary.Range("C11").Value = 680
$wsSummary.Range("D11").Value = 648
$wsSummary.Range("F11").Value = -1.24
$wsSummary.Range("I11").Value = -0.8
$wsSummary.Range("L11").Value = 0.79
$wsSummary.Range("C12").Value = 714
$wsSummary.Range("D12").Value = 676
$wsSummary.Range("E12").Value = 63.76
$wsSummary.Range("F12").Value = -0.74
$wsSummary.Range("H12").Value = 1
$wsSummary.Range("I12").Value = 0.1
$wsSummary.Range("J12").Value = 0.34
$wsSummary.Range("K12").Value = 0.06
$wsSummary.Range("L12").Value = 1.17
$wsSummary.Range("C13").Value = 714
$wsSummary.Range("D13").Value = 643
$wsSummary.Range("E13").Value = 55.37
$wsSummary.Range("F13").Value = -0.29
$wsSummary.Range("I13").Value = 0.27
$wsSummary.Range("J13").Value = 0.38
$wsSummary.Range("L13").Value = 1.29
$wsSummary.Range("C14").Value = 714
$wsSummary.Range("D14").Value = 623
$wsSummary.Range("E14").Value = 51.04
$wsSummary.Range("I14").Value = 0.45
$wsSummary.Range("J14").Value = 0.32
$wsSummary.Range("L14").Value = 1.09
$wsSummary.Range("C15").Value = 714
$wsSummary.Range("E15").Value = 47.62
$wsSummary.Range("F15").Value = 0.11
$wsSummary.Range("H15").Value = 1.04
$wsSummary.Range("I15").Value = 0.96
$wsSummary.Range("J15").Value = 0.32
$wsSummary.Range("L15").Value = 1.02
$wsSummary.Range("C16").Value = 714
$wsSummary.Range("E16").Value = 42.47
$wsSummary.Range("F16").Value = 0.35
$wsSummary.Range("H16").Value = 1.04
$wsSummary.Range("I16").Value = 0.96
$wsSummary.Range("J16").Value = 0.33
$wsSummary.Range("K16").Value = 0.03
$wsSummary.Range("L16").Value = 1.02
$wsSummary.Range("C17").Value = 714
$wsSummary.Range("D17").Value = 336
$wsSummary.Range("F17").Value = -0.26
$wsSummary.Range("H17").Value = 0.92
$wsSummary.Range("I17").Value = -1.16
$wsSummary.Range("J17").Value = 0.56
$wsSummary.Range("K17").Value = 0.1
$wsSummary.Range("L17").Value = 0.88
$wsSummary.Range("C18").Value = 706
$wsSummary.Range("D18").Value = 671
$wsSummary.Range("E18").Value = 39.64
$wsSummary.Range("F18").Value = 0.8
$wsSummary.Range("I18").Value = 0.65
$wsSummary.Range("J18").Value = 0.32
$wsSummary.Range("L18").Value = 1.11
$wsSummary.Range("C19").Value = 706
$wsSummary.Range("D19").Value = 678
$wsSummary.Range("E19").Value = 37.46
$wsSummary.Range("F19").Value = 0.92
$wsSummary.Range("I19").Value = 1.32
$wsSummary.Range("J19").Value = 0.28
$wsSummary.Range("K19").Value = 0.07
$wsSummary.Range("L19").Value = 0.97
$wsSummary.Range("C20").Value = 706
$wsSummary.Range("D20").Value = 678
$wsSummary.Range("H20").Value = 0.95
$wsSummary.Range("I20").Value = -1.49
$wsSummary.Range("J20").Value = 0.36
$wsSummary.Range("L20").Value = 0.86
$wsSummary.Range("C21").Value = 706
$wsSummary.Range("D21").Value = 676
$wsSummary.Range("E21").Value = 26.48
$wsSummary.Range("F21").Value = 1.56
$wsSummary.Range("H21").Value = 0.94
$wsSummary.Range("I21").Value = -1.17
$wsSummary.Range("L21").Value = 1.71
$wsSummary.Range("C22").Value = 706
$wsSummary.Range("D22").Value = 676
$wsSummary.Range("E22").Value = 23.08
$wsSummary.Range("F22").Value = 1.8
$wsSummary.Range("L22").Value = 1.41

# --- "model_fit" sheet (sheet2) numeric corrections ---
$wsModelFit = $wb.Worksheets.Item("model_fit")
$wsModelFit.Range("D2").Value = 24489
$wsModelFit.Range("E2").Value = 24547
$wsModelFit.Range("F2").Value = 24711
$wsModelFit.Range("G2").Value = 0.666
$wsModelFit.Range("H2").Value = 0.569
$wsModelFit.Range("D3").Value = 24385
$wsModelFit.Range("E3").Value = 24483
$wsModelFit.Range("F3").Value = 24760
$wsModelFit.Range("G3").Value = 0.679
$wsModelFit.Range("H3").Value = 0.576

# --- "steps" sheet (sheet3) shared-string corrections ---
$wsSteps = $wb.Worksheets.Item("steps")

# Plain text updates (values are not numeric-looking, so they stay text automatically)
$wsSteps.Range("B3").Value = "0.75 (0.062)"
$wsSteps.Range("B4").Value = "0.92 (0.092)"
$wsSteps.Range("B5").Value = "1.39 (0.122)"
$wsSteps.Range("C5").Value = "-1.03 (0.142)"

# Updates where the new text looks like a plain number; force text storage
# (format as Text, assign the value, then restore the default "Normal" style
# so no stray number-format is left applied to the cell)
$wsSteps.Range("C3").NumberFormat = "@"
$wsSteps.Range("C3").Value = "-0.75"
$wsSteps.Range("C3").Style = "Normal"

$wsSteps.Range("D4").NumberFormat = "@"
$wsSteps.Range("D4").Value = "0.09"
$wsSteps.Range("D4").Style = "Normal"

$wsSteps.Range("D5").NumberFormat = "@"
$wsSteps.Range("D5").Value = "-0.36"
$wsSteps.Range("D5").Style = "Normal"
